# Insert a new data row at row 311 (pushing the existing rows 311-370 down
# to 312-371) and populate it with a new weekly price observation for
# "Vega Monumental Concepción" / Zanahoria.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 311..370 down one row, creating a blank row 311.
$ws.Rows.Item(311).Insert()

# Populate the newly inserted row 311 with the new record.
$ws.Range("A311").Value = 11
$ws.Range("B311").Value = "Vega Monumental Concepción"
$ws.Range("C311").Value = "Bíobío"
$ws.Range("D311").Value = 45015
$ws.Range("E311").Value = 8
$ws.Range("F311").Value = 100114013
$ws.Range("G311").Value = "Zanahoria"
$ws.Range("H311").Value = "Sin especificar"
$ws.Range("I311").Value = "Primera"
$ws.Range("J311").Value = 220
$ws.Range("K311").Value = 4000
$ws.Range("L311").Value = 4500
$ws.Range("M311").Value = 4273
$ws.Range("N311").Value = "$/saco 20 kilos"
$ws.Range("O311").Value = "Región de Ñuble"
$ws.Range("P311").Value = 214
$ws.Range("Q311").Value = 20
$ws.Range("R311").Value = "Hortaliza"
